$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.133.79'
$ws.Range("E2").Value = '  -0.26%  '

$ws.Range("D3").Value = '3.132.81'
$ws.Range("E3").Value = '  -0.76%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.07'
$ws.Range("E5").Value = '  +0.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.03'
$ws.Range("E6").Value = '  -2.83%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").Value = '3.129.79'
$ws.Range("E8").Value = '  -0.70%  '

$ws.Range("E9").Value = '  +0.51%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.150'
$ws.Range("E10").Value = '  -0.18%  '

$ws.Range("E11").Value = '  -2.40%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.476'
$ws.Range("E12").Value = '  +0.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000256'
$ws.Range("E13").Value = '  +2.26%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.46'
$ws.Range("E14").Value = '  -1.86%  '

$ws.Range("D15").Value = '3.651.08'
$ws.Range("E15").Value = '  -0.65%  '

$ws.Range("E16").Value = '  +2.64%  '

$ws.Range("D17").Value = '64.134.21'
$ws.Range("E17").Value = '  -0.29%  '

$ws.Range("D18").Value = '3.126.41'
$ws.Range("E18").Value = '  -0.92%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.86'
$ws.Range("E19").Value = '  -1.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '476.34'
$ws.Range("E20").Value = '  -0.66%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.67'
$ws.Range("E21").Value = '  -0.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.721'
$ws.Range("E22").Value = '  +1.61%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.80'
$ws.Range("E23").Value = '  +1.50%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.61'
$ws.Range("E24").Value = '  -0.76%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.46'
$ws.Range("E25").Value = '  +2.21%  '

$ws.Range("E26").Value = '  -0.02%  '

$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.78'
$ws.Range("E27").Value = '  -3.22%  '

$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.54'
$ws.Range("E28").Value = '  +1.41%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.42'
$ws.Range("E29").Value = '  +8.66%  '

$ws.Range("E30").Value = '  +2.99%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.08'
$ws.Range("E31").Value = '  -5.02%  '

$ws.Range("E32").Value = '  -0.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.61'
$ws.Range("E33").Value = '  +1.59%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.62'
$ws.Range("E34").Value = '  -4.42%  '

$ws.Range("E35").Value = '  +0.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.93'
$ws.Range("E36").Value = '  -0.75%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.38'
$ws.Range("E37").Value = '  -3.42%  '

$ws.Range("D38").Value = '0.0₃0742'
$ws.Range("E38").Value = '  +4.11%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '453.02'
$ws.Range("E39").Value = '  +0.59%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.99'
$ws.Range("E40").Value = '  +2.99%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0395'
$ws.Range("E41").Value = '  -0.09%  '

$ws.Range("E42").Value = '  -0.78%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.33'
$ws.Range("E43").Value = '  -1.18%  '

$ws.Range("D44").Value = '2.867.62'
$ws.Range("E44").Value = '  +0.97%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.264'
$ws.Range("E45").Value = '  -1.14%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.26'
$ws.Range("E46").Value = '  -0.74%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.43'
$ws.Range("E47").Value = '  +4.79%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.40'
$ws.Range("E48").Value = '  +0.09%  '

$ws.Range("E49").Value = '  +0.07%  '

$ws.Range("E50").Value = '  +0.03%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.75'
$ws.Range("E51").Value = '  +2.13%  '
